$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = 45235
$ws.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B39").Value = "11:20"
$ws.Range("C39").Value = 1819
$ws.Range("D39").Value = "amazon"
$ws.Range("E39").Value = "preto"
